# Add Q3-2022 data to the "300428-立中集团" workbook.
#
# Plan:
#  1. "总计" (sheet 1) summary sheet: insert a new row right below the
#     header for the "2022-Q3" totals, shifting the existing quarters
#     down by one row (and renumbering their running index in column A).
#  2. Duplicate the existing "2022-Q2" detail sheet (so the brand new
#     "2022-Q3" sheet starts out with identical formatting/styles), put
#     the duplicate back in the "2022-Q2" slot, and turn the original
#     sheet object into "2022-Q3" with the new fund-holding rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "总计" overview sheet - insert the new 2022-Q3 row at the top
#    of the data (row 2), pushing everything else down.
# ---------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$wsTotal.Rows.Item(2).Insert()

# Re-apply the bold/centered/bordered "index" style (column A) that the
# other data rows use - copy it down from the row that used to be row 2
# (now row 3) onto the freshly inserted row 2.
$wsTotal.Cells.Item(3, 1).Copy()
$wsTotal.Cells.Item(2, 1).PasteSpecial(-4122)
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 3
$wsTotal.Cells.Item(2, 4).Value = 0.03

# Renumber the running index (column A) of the quarters that shifted down.
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(4, 1).Value = 2
$wsTotal.Cells.Item(5, 1).Value = 3

# ---------------------------------------------------------------
# 2) Detail sheets - insert "2022-Q3" before "2022-Q2".
# ---------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item(2)

# Duplicate "2022-Q2" right after itself; the copy will keep the
# original quarter's data/name and slot back in as "2022-Q2".
$wsQ2.Copy($null, $wsQ2)

$wsNewQ3 = $wb.Worksheets.Item(2)
$wsQ2Copy = $wb.Worksheets.Item(3)

$wsNewQ3.Name = "2022-Q3"
$wsQ2Copy.Name = "2022-Q2"

# Wipe the (copied) 2022-Q2 rows from the sheet that is becoming 2022-Q3,
# keeping only the header row, then write the new fund-holding data.
$wsNewQ3.Range("A5:H23").Clear()

# Force columns B:G to be stored as text (fund codes/percentages keep
# leading zeros / exact decimal formatting, matching the source data).
$wsNewQ3.Range("B2:G4").NumberFormat = "@"

$wsNewQ3.Cells.Item(2, 1).Value = 0
$wsNewQ3.Cells.Item(2, 2).Value = "011765"
$wsNewQ3.Cells.Item(2, 3).Value = "兴银高端制造混合A"
$wsNewQ3.Cells.Item(2, 4).Value = "0.57"
$wsNewQ3.Cells.Item(2, 5).Value = "92.99"
$wsNewQ3.Cells.Item(2, 6).Value = "2.78"
$wsNewQ3.Cells.Item(2, 7).Value = "0.0158"
$wsNewQ3.Cells.Item(2, 8).Value = 9

$wsNewQ3.Cells.Item(3, 1).Value = 1
$wsNewQ3.Cells.Item(3, 2).Value = "011766"
$wsNewQ3.Cells.Item(3, 3).Value = "兴银高端制造混合C"
$wsNewQ3.Cells.Item(3, 4).Value = "0.34"
$wsNewQ3.Cells.Item(3, 5).Value = "92.99"
$wsNewQ3.Cells.Item(3, 6).Value = "2.78"
$wsNewQ3.Cells.Item(3, 7).Value = "0.0095"
$wsNewQ3.Cells.Item(3, 8).Value = 9

$wsNewQ3.Cells.Item(4, 1).Value = 2
$wsNewQ3.Cells.Item(4, 2).Value = "005146"
$wsNewQ3.Cells.Item(4, 3).Value = "兴银丰润灵活配置混合"
$wsNewQ3.Cells.Item(4, 4).Value = "0.04"
$wsNewQ3.Cells.Item(4, 5).Value = "92.81"
$wsNewQ3.Cells.Item(4, 6).Value = "3.14"
$wsNewQ3.Cells.Item(4, 7).Value = "0.0013"
$wsNewQ3.Cells.Item(4, 8).Value = 10

# Cells were written as General first (via NumberFormat forcing text
# entry) - reset the visual style back to Normal so they match the
# un-styled data cells used throughout the workbook, while keeping the
# values stored as text.
$wsNewQ3.Range("B2:G4").Style = "Normal"

Write-Output "2022-Q3 sheet added"
